# Optimized NJ algo and nice print visuals
#
# The sheet's genotype table (rows 19-35) is restructured: a new sample
# "weird1" is inserted at the top of that block, the existing "350xyz",
# 14196 and 14206 sample blocks shift down, and a new "39copy" sample
# block is appended at the end. Net effect: rows 19-35 (17 rows) become
# rows 19-47 (29 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old rows 19-35 data block outright; we rebuild rows 19-47 from
# scratch below so there's no ambiguity about which old row maps to which
# new row.
$ws.Range("A19:C35").ClearContents() | Out-Null

# NOTE on shared-string ordering: the workbook's shared string table is
# rebuilt from the order in which distinct string values are first
# assigned to cells. The target file needs "39copy" to land at index 9
# and "weird1" at index 10, so we must write a "39copy" cell before any
# "weird1" cell, even though "weird1" appears in an earlier row on the
# sheet. We do that by populating row 42 (the first "39copy" row) first.
$newRows = [ordered]@{
    42 = @("39copy", "Tr1",  500)

    19 = @("weird1", "Tr1",  194)
    20 = @("weird1", "Tr15", 150)
    21 = @("weird1", "Tr17", 92)
    22 = @("weird1", "Tr5",  651)
    23 = @("weird1", "Tr5",  651)
    24 = @("weird1", "Tr5",  795)

    25 = @("350xyz", "Tr1",  129)
    26 = @("350xyz", "Tr15", 150)
    27 = @("350xyz", "Tr17", 98)
    28 = @("350xyz", "Tr17", 92)
    29 = @("350xyz", "Tr17", 105)
    30 = @("350xyz", "Tr2",  808)
    31 = @("350xyz", "Tr5",  795)

    32 = @(14196, "Tr1",  490)
    33 = @(14196, "Tr15", 150)
    34 = @(14196, "Tr17", 98)
    35 = @(14196, "Tr2",  687)
    36 = @(14196, "Tr5",  795)

    37 = @(14206, "Tr1",  194)
    38 = @(14206, "Tr15", 150)
    39 = @(14206, "Tr17", 92)
    40 = @(14206, "Tr2",  627)
    41 = @(14206, "Tr5",  651)

    43 = @("39copy", "Tr15", 142)
    44 = @("39copy", "Tr17", 92)
    45 = @("39copy", "Tr17", 92)
    46 = @("39copy", "Tr2",  808)
    47 = @("39copy", "Tr5",  651)
}

foreach ($r in $newRows.Keys) {
    $vals = $newRows[$r]
    $ws.Range("A$r").Value = $vals[0]
    $ws.Range("B$r").Value = $vals[1]
    $ws.Range("C$r").Value = $vals[2]
}

# Workbook/sheet VBA code names, as reflected in the target workbookPr /
# sheetPr elements.
$wb.CodeName = "ThisWorkbook"
$ws.CodeName = "Sheet1"

# Match the new selected cell shown in the target sheetView.
$ws.Range("L9").Select() | Out-Null
